# The deck's embedded "Integral" master theme (ppt/theme/theme1.xml) is
# being swapped for a stock "Office Theme" colour palette (the palette
# that used to live alongside it as the Notes Master's theme,
# ppt/theme/theme2.xml). Re-colour the 12 theme colour slots (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) through the slide's
# ThemeColorScheme - this writes straight into the <a:clrScheme> of the
# master theme without disturbing anything else (names/fonts/format
# scheme already agree between the two themes in this deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colours, in clrScheme order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. PowerPoint RGB() values are packed as
# 0x00BBGGRR, i.e. the reverse byte order of the familiar #RRGGBB.
$tcs.Colors(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1      000000
$tcs.Colors(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1      FFFFFF
$tcs.Colors(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2      44546A
$tcs.Colors(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2      E7E6E6
$tcs.Colors(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1  5B9BD5
$tcs.Colors(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2  ED7D31
$tcs.Colors(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3  A5A5A5
$tcs.Colors(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4  FFC000
$tcs.Colors(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5  4472C4
$tcs.Colors(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6  70AD47
$tcs.Colors(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink    0563C1
$tcs.Colors(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink 954F72
